$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# Header date
Replace-Text "2023-11-29 Wednesday" "2023-11-30 Thursday"

# Table cell problems (order matters: 99÷8= must be replaced before 31÷5= becomes 99÷8=)
Replace-Text "77÷4=" "35÷7="
Replace-Text "78÷6=" "97÷8="
Replace-Text "29÷9=" "11÷6="
Replace-Text "64÷6=" "75÷9="
Replace-Text "50÷9=" "34÷5="
Replace-Text "21÷4=" "59÷4="
Replace-Text "82÷6=" "47÷6="
Replace-Text "79÷3=" "78÷3="
Replace-Text "36÷9=" "22÷2="
Replace-Text "67÷4=" "84÷8="
Replace-Text "29÷4=" "51÷2="
Replace-Text "98÷7=" "90÷8="
Replace-Text "60÷2=" "69÷2="
Replace-Text "44÷8=" "35÷3="
Replace-Text "99÷8=" "59÷2="
Replace-Text "46÷9=" "13÷7="
Replace-Text "46÷8=" "71÷8="
Replace-Text "82÷9=" "67÷3="
Replace-Text "15÷6=" "17÷9="
Replace-Text "99÷2=" "63÷2="
Replace-Text "40÷6=" "32÷2="
Replace-Text "31÷5=" "99÷8="
Replace-Text "83÷9=" "92÷6="
Replace-Text "38÷9=" "51÷6="
Replace-Text "93÷3=" "18÷7="
